$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank spacer row that follows the "Performance Requirements"
# section (old row 7), so the trailing blank row moves up to become row 7.
$ws.Rows("7").Delete()

# Fill in the new Performance Requirements bullet text in what is now row 6
# (the blank placeholder row right after "Performance Requirements").
$ws.Range("A6").Value = "The game actions and updates SHOULD have fast response times."

# Give it its own indented (non-bold) style, distinct from the bold section
# header style: Aptos 12, left/center aligned, indent level 2.
$ws.Range("A6").Font.Name = "Aptos"
$ws.Range("A6").Font.Size = 12
$ws.Range("A6").Font.Bold = $false
$ws.Range("A6").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A6").VerticalAlignment = -4108    # xlCenter
$ws.Range("A6").IndentLevel = 2

# Update the selection to match the saved state of the workbook.
$ws.Range("C16").Select() | Out-Null
